# Actualización automática hashcode
# This script updates the "hashcode" values (column B) in the active worksheet
# to match the new values produced by the automated hashcode regeneration run.
# Each entry below maps a worksheet cell to its expected previous value and the
# new value that should replace it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "B11"; Old = "8d34570a89896d1e7487e15264d4430d"; New = "27727890b1e1b6a6913a83c5b04b29a4" }
    @{ Cell = "B29"; Old = "c9ab33bf5bace551342158f57f5fb0c5"; New = "283a43421d3d619311cfd0592b2dd6eb" }
    @{ Cell = "B121"; Old = "27ce3918723a74c22be7d3b4776af7d0"; New = "a75a4d9af7223344b490c2aca5cdac25" }
    @{ Cell = "B126"; Old = "30992a194a56e3775d7bc9fa5a64bc24"; New = "51bbf56d85cc17f3c8cb856bf4fd262d" }
    @{ Cell = "B136"; Old = "145f6cdd9e574970a49058607a4c57c6"; New = "5e3fe43d9be5b777179b6c69eea2d63f" }
    @{ Cell = "B159"; Old = "dbfc21f7e94c2499a7e91e097f364003"; New = "17e6f09fd8ea8a8972bc475df817080f" }
    @{ Cell = "B169"; Old = "d8e2d3b430620fbcc36650018a5d213d"; New = "6afcb86346c0f16cac73003425cae14d" }
    @{ Cell = "B191"; Old = "abec732590cbb485771a4e1b3fca3502"; New = "3a13e7d435e81d1a9016877eee3af917" }
    @{ Cell = "B227"; Old = "366679d9cd102f7c634ebffd2d642faa"; New = "79d7ac27c02b8ee4b146a8ebaf9cdac1" }
    @{ Cell = "B232"; Old = "2ad3ae0d1889ca9238638c3c5377ba7a"; New = "ae22bcdb5a3d16e8e1bb7667b80435a8" }
    @{ Cell = "B281"; Old = "7f6ab24a2600337270ff3e0396ae3efd"; New = "91d6cecafdef3ad37838abc58fd1f3c8" }
    @{ Cell = "B293"; Old = "66fae7c05456a4b684f7c16d5b50be85"; New = "8cb4f938f3e6a3f50370cb099b1625d5" }
    @{ Cell = "B302"; Old = "0f1ef506e706195dbd93c49065f789b1"; New = "d263c9cd625e0cc36308d3fec4350e23" }
    @{ Cell = "B339"; Old = "4355b8ccd9f3d91560badc347230afcd"; New = "1e506b1f2a033ed20095cbdd53afc20a" }
    @{ Cell = "B419"; Old = "2ee5add6736bc97726d8045230c25adb"; New = "afba4ee92bb44bede48ddf483ac24705" }
    @{ Cell = "B460"; Old = "ef3bb11c9a11290215fab20c3653025e"; New = "0cd8625297c32aba25b0f61545f1b53e" }
    @{ Cell = "B478"; Old = "0e421a028fe726870a018a31b7132a98"; New = "19b25a4ce25f6f97839a85d363ab88b0" }
    @{ Cell = "B480"; Old = "18e444a0140e6b442fd1939ef7a91505"; New = "54047bec7956934d2f51b05c58bf2b32" }
    @{ Cell = "B500"; Old = "90638a5840cb2ea45547ac598d99705e"; New = "59328d6fbee2ac587678815c09af1874" }
    @{ Cell = "B501"; Old = "10add39a694426657601535a2ecb2c04"; New = "2f3dfc70d7f041da9765e62f76ca913a" }
    @{ Cell = "B502"; Old = "2be9b891f2e904e9681becfeefd8ad95"; New = "81629ac93065ab0b8af54b4a0aeeec72" }
    @{ Cell = "B517"; Old = "d58681c86cbed19c395aab18d70338ab"; New = "4411e56c2ff7e6ec8787d8f6be166e8b" }
    @{ Cell = "B550"; Old = "8aab137630c87b0adee966d8555f7e13"; New = "345984d1f1a72d556b2fb2538b0e94aa" }
    @{ Cell = "B572"; Old = "2829c5fc1f67e224165dc8d654e289f4"; New = "0751fcd52a01e68b0dea88477cc78546" }
    @{ Cell = "B616"; Old = "078638d89707ef761041c1aa1f6eb798"; New = "cf51451dd6f5b3073cd680b0a9c8f098" }
    @{ Cell = "B627"; Old = "0225aa8685f6b6a513936ce0d53587e9"; New = "cd0f810a0814b71df06adc86d49f9165" }
    @{ Cell = "B629"; Old = "00d68d50c3de3d47c92bdab22d9dc903"; New = "8e135f17d024197e1fee484b3eb87bd1" }
    @{ Cell = "B649"; Old = "759613b2f4e599e5bbf90a4d43e40cc9"; New = "3e72d49f2d5a1c266973b510c1bc866a" }
    @{ Cell = "B655"; Old = "6a5e3c6b8da31df5f747f3f32e2ebcf8"; New = "d6d55401dea2dc036bcb028447293785" }
    @{ Cell = "B665"; Old = "8d7ab49717672ca1dce808920279e201"; New = "1ba24c61578dfbe6dd75691af4a3de32" }
    @{ Cell = "B819"; Old = "ddcecae74f700d34aeb688e4eafe9966"; New = "f918429f8f38492013789bfd11f54108" }
    @{ Cell = "B830"; Old = "878f501c6fcfbb24100b756563e49341"; New = "39131b3cfdad3487567b097fc174ea20" }
    @{ Cell = "B835"; Old = "44a1dc031076aedec8ddf2465a2c79d5"; New = "6c0c01f6b02ef111a430a37b418b5556" }
    @{ Cell = "B862"; Old = "15adcc8626573003a2667afe259f8d2e"; New = "56ad9242b497ae392e8130d0697a5abd" }
    @{ Cell = "B874"; Old = "c9c849f03081bb7a17b5eba5feebb7ea"; New = "d878f735a89572d2273c1e98708e28dd" }
)

$updated = 0
$mismatched = 0

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $current = $cell.Value2

    if ($current -ne $change.Old) {
        $mismatched++
        Write-Host "Warning: $($change.Cell) expected '$($change.Old)' but found '$current'"
    }

    $cell.Value2 = $change.New
    $updated++
}

Write-Host "Updated $updated cell(s); $mismatched mismatch(es) against expected prior value."
